$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new rows of data (row 29 and row 30)
$ws.Range("A29").Value = "2019年4月25日23:23:35"
$ws.Range("B29").Value = "周四"
$ws.Range("C29").Value = "业务实体domain Entity ER建模"
$ws.Range("D29").Value = "10:30--12:00"

$ws.Range("C30").Value = "star项目构建（base，dao）Entity and Test"
$ws.Range("D30").Value = "18:30--21:10"

# Update the selection to match the recorded cursor position
$ws.Range("C30").Select()
